$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Resumen": update the zone label and the associated maximum metric.
# ---------------------------------------------------------------------------
$wsResumen = $wb.Worksheets.Item("Resumen")
$wsResumen.Range("B2").Value = "Z2"
$wsResumen.Range("C2").Value = 514.402816520705

# ---------------------------------------------------------------------------
# Sheet "Solucion": the randomized constructive method produced a new
# Pedido -> Salida assignment order. Overwrite rows 2..41 (A:B) in place.
# ---------------------------------------------------------------------------
$wsSolucion = $wb.Worksheets.Item("Solucion")

$solucionData = @(
    @("Pedido_14", "S001"),
    @("Pedido_26", "S021"),
    @("Pedido_13", "S031"),
    @("Pedido_28", "S011"),
    @("Pedido_35", "S022"),
    @("Pedido_37", "S002"),
    @("Pedido_7", "S012"),
    @("Pedido_20", "S032"),
    @("Pedido_5", "S003"),
    @("Pedido_10", "S023"),
    @("Pedido_16", "S033"),
    @("Pedido_36", "S013"),
    @("Pedido_17", "S004"),
    @("Pedido_22", "S024"),
    @("Pedido_40", "S014"),
    @("Pedido_15", "S034"),
    @("Pedido_4", "S005"),
    @("Pedido_27", "S025"),
    @("Pedido_19", "S015"),
    @("Pedido_2", "S035"),
    @("Pedido_8", "S006"),
    @("Pedido_21", "S016"),
    @("Pedido_38", "S007"),
    @("Pedido_30", "S026"),
    @("Pedido_12", "S036"),
    @("Pedido_1", "S017"),
    @("Pedido_33", "S027"),
    @("Pedido_25", "S008"),
    @("Pedido_29", "S037"),
    @("Pedido_32", "S018"),
    @("Pedido_31", "S028"),
    @("Pedido_6", "S009"),
    @("Pedido_23", "S038"),
    @("Pedido_9", "S029"),
    @("Pedido_3", "S019"),
    @("Pedido_34", "S039"),
    @("Pedido_18", "S030"),
    @("Pedido_11", "S010"),
    @("Pedido_39", "S020"),
    @("Pedido_24", "S040")
)

for ($i = 0; $i -lt $solucionData.Length; $i++) {
    $row = $i + 2
    $wsSolucion.Cells.Item($row, 1).Value = $solucionData[$i][0]
    $wsSolucion.Cells.Item($row, 2).Value = $solucionData[$i][1]
}

# ---------------------------------------------------------------------------
# Sheet "Metricas": refresh the per-zone timing metrics.
# ---------------------------------------------------------------------------
$wsMetricas = $wb.Worksheets.Item("Metricas")
$wsMetricas.Range("B2").Value = 514.2673802573253
$wsMetricas.Range("B3").Value = 514.402816520705
